# FanPowerAllowances-T24N_2022.xlsx edit
# Adds 4 new columns (Exhaust Systems Base Allowance / Supply ERV / Return ERV / Return Filter)
# to the "TABLE T24N_2022FanPwrIdxAdj" block, pushing the old "SZVAV" column from K to N.
# The second table below (T24N_2022BaseFanPwrIdx, rows 13-19) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Move the existing "SZVAV" column (K3:K10) to its new home at column N,
#    carrying values + number formatting with it, before we overwrite K:M.
$ws.Range("K3:K10").Copy($ws.Range("N3"))

# 2. New header row 3 text (row also becomes slightly shorter - 90pt instead of 105pt)
$ws.Range("J3").Value = "Exhaust Systems Base Allowance"
$ws.Range("K3").Value = "Supply - Energy Recovery (Enthalpy Recovery Ratio ? 0.60 and <0.65)  "
$ws.Range("L3").Value = "Return - Energy Recovery (Enthalpy Recovery Ratio ? 0.60 and <0.65)  "
$ws.Range("M3").Value = "Return - Filter (any MERV value)"
$ws.Rows(3).RowHeight = 90

# 3. New header row 4 (internal field names) - row grows to 30pt and gets wrap text
$ws.Range("J4").Value = "ExhBaseAllow"
$ws.Range("K4").Value = "ERVSupply"
$ws.Range("L4").Value = "ERVReturn"
$ws.Range("M4").Value = "RetFilter"
$ws.Range("J4:M4").WrapText = $true
$ws.Rows(4).RowHeight = 30

# 4. New data values for rows 5-10 (plain numbers, no special number format)
$ws.Range("J5").Value = 0.221
$ws.Range("K5").Value = 0.184
$ws.Range("L5").Value = 0.19
$ws.Range("M5").Value = 0.046

$ws.Range("J6").Value = 0.246
$ws.Range("K6").Value = 0.155
$ws.Range("L6").Value = 0.163
$ws.Range("M6").Value = 0.041

$ws.Range("J7").Value = 0.236
$ws.Range("K7").Value = 0.144
$ws.Range("L7").Value = 0.146
$ws.Range("M7").Value = 0.036

$ws.Range("J8").Value = 0.186
$ws.Range("K8").Value = 0.19
$ws.Range("L8").Value = 0.191
$ws.Range("M8").Value = 0.046

$ws.Range("J9").Value = 0.184
$ws.Range("K9").Value = 0.163
$ws.Range("L9").Value = 0.166
$ws.Range("M9").Value = 0.041

$ws.Range("J10").Value = 0.19
$ws.Range("K10").Value = 0.146
$ws.Range("L10").Value = 0.148
$ws.Range("M10").Value = 0.036

# 5. Column widths: D:I stay at 20 chars; J:N become the new, wider 22-char columns
$ws.Range("J1:N1").EntireColumn.ColumnWidth = 22

# 6. Selection ends on M15, matching the saved file
$ws.Range("M15").Select()
